$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "CreateRecipient" (tab 1)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("CreateRecipient")
$ws1.Rows.Item(1).ClearFormats()
$ws1.Rows.Item(2).ClearFormats()
$ws1.Range("A1").Value = "Name1"
$ws1.Range("B1").Value = "Name2"
$ws1.Range("A2").Value = "'Fax Address"
$ws1.Range("B2").Value = "'Recipient"
$ws1.Range("C2").Value = "'9987288"
$ws1.Range("A1").Select()

# ---------------------------------------------------------------------------
# Sheet "EditRecipient" (tab 2)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("EditRecipient")
$ws2.Rows.Item(1).ClearFormats()
$ws2.Rows.Item(2).ClearFormats()
$ws2.Range("A1").Value = "Name1"
$ws2.Range("B1").Value = "Name2"
$ws2.Range("A2").Value = "'Fax Address"
$ws2.Range("B2").Value = "'Recipient"
$ws2.Range("C2").Value = "'9987288"
$ws2.Range("D2").Value = "'Fax Address Updated"
$ws2.Range("E2").Value = "Modified"
$ws2.Range("D3").Select()

# ---------------------------------------------------------------------------
# Sheet "DeleteRecipient" (tab 3)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("DeleteRecipient")
$ws3.Rows.Item(1).ClearFormats()
$ws3.Rows.Item(2).ClearFormats()
$ws3.Range("A1").Value = "Name1"
$ws3.Range("B1").Value = "Name2"
$ws3.Range("A2").Value = "'Fax Address Updated"
$ws3.Range("B2").Value = "'Recipient"
$ws3.Range("C2").Value = "'9987288"
$ws3.Range("D2").Value = "'Deleted"
$ws3.Range("B1").Select()

# ---------------------------------------------------------------------------
# Sheet "AddressCreate" (tab 4) - becomes the active tab
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("AddressCreate")
# D2 carries style index 2 (quotePrefix font) - keep it with a leading quote
$ws4.Range("D2").Value = "'sample2 delete<12345>"
$ws4.Range("E2").Value = "'Fax Address Recipient<9987288>,sample2 delete<12345>"
$ws4.Columns.Item(3).ColumnWidth = 29.92
$ws4.Columns.Item(4).ColumnWidth = 20.75
$ws4.Activate()
$ws4.Range("E6").Select()

# ---------------------------------------------------------------------------
# Sheet "EditAddressBook" (tab 5)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("EditAddressBook")
# C2/E2 carry style index 1 - keep it with a leading quote
$ws5.Range("C2").Value = "'Fax Address Recipient<9987288>,SampleData<123467>"
$ws5.Range("E2").Value = "'Fax Address Updated Recipient<9987288>,SampleData<123467>"
$ws5.Range("C2").Select()

# ---------------------------------------------------------------------------
# Sheet "DeleteAddressBook" (tab 6)
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("DeleteAddressBook")
# A2 has no style, but C2 carries style index 1 - keep it with a leading quote
$ws6.Range("A2").Value = "Address Book Updated"
$ws6.Range("C2").Value = "'Fax Address Updated Recipient<9987288>,SampleData<123467>"
$ws6.Columns.Item(1).ColumnWidth = 21.42
$ws6.Columns.Item(3).ColumnWidth = 49.59
$ws6.Range("C7").Select()

# Make sure AddressCreate ends up the active sheet/tab (activeTab=3 in workbook.xml)
$ws4.Activate()
